# PR "Kho nhận" column: insert a new column C ("Warehouse received") into the
# purchase-request import template, between "Lệnh sản xuất" (B) and
# "Trung tâm chi phí" (C, which shifts to D along with everything after it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting shifts C:O -> D:P to the right and copies column B's formatting
# (incl. width) into the freshly-opened column C.
$ws.Columns("C:C").Insert()

# The insert also clones row 4's lone cell (B4) into the new column; the
# template only ever had a cell in B4, so drop the stray C4 that appeared.
$ws.Range("C4").Clear()

# New data row value and header label for the inserted column.
$ws.Range("C2").Value = 1075
$ws.Range("C1").Value = "Kho nhận"

# Both cells are formatted as Text, matching the rest of the template's
# "code-like" columns (e.g. the order number in column D/"Mã vụ việc" style).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C1").NumberFormat = "@"

# Match the author's final selection/cursor position.
$ws.Range("G14").Select() | Out-Null
